# DU-0001-234987347.xlsx fixture cleanup:
#  - rename the first sheet from the default LibreOffice "Folha1" to the
#    descriptive "Concentrations"
#  - update the remembered cell selections: sheet1 ends up with G44 selected
#    (instead of the stale B3), while sheet2/sheet3 just keep A1 selected
#    (instead of the stale B3 remembered alongside A1)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Concentrations"

# Touch sheet2 / sheet3 so their saved selection drops the old "B3" memory
# and is left pointing only at A1.
[void]$ws2.Activate()
[void]$ws2.Range("A1").Select()

[void]$ws3.Activate()
[void]$ws3.Range("A1").Select()

# Finish on sheet1 with G44 selected/active - this is the sheet left visible
# when the workbook is reopened.
[void]$ws1.Activate()
[void]$ws1.Range("G44").Select()
